$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrB = New-Object 'object[,]' 24,1
$arrB[0,0] = 0.1424986629922103
$arrB[1,0] = 0.1329151440306902
$arrB[2,0] = 0.1270977887459424
$arrB[3,0] = 0.1247442048677527
$arrB[4,0] = 0.1243544287756464
$arrB[5,0] = 0.1270659782516219
$arrB[6,0] = 0.1391804819419917
$arrB[7,0] = 0.1634609971099081
$arrB[8,0] = 0.181611468663931
$arrB[9,0] = 0.1899347012111008
$arrB[10,0] = 0.193095882340927
$arrB[11,0] = 0.1924146527230732
$arrB[12,0] = 0.1901945872134405
$arrB[13,0] = 0.1888359449413173
$arrB[14,0] = 0.1810688460488024
$arrB[15,0] = 0.1763208645652838
$arrB[16,0] = 0.1735962221262071
$arrB[17,0] = 0.1726747883934081
$arrB[18,0] = 0.1768256477894568
$arrB[19,0] = 0.1908464220682902
$arrB[20,0] = 0.2000642324690318
$arrB[21,0] = 0.1951396023207934
$arrB[22,0] = 0.1765974195935485
$arrB[23,0] = 0.1568371918552316
$ws.Range("B2:B25").Value = $arrB

$arrDH = New-Object 'object[,]' 24,5
$arrDH[0,0] = 0.01099854060988292
$arrDH[0,1] = 0.07835222606900327
$arrDH[0,2] = 1.063950388120105
$arrDH[0,3] = 1.050397112036137
$arrDH[0,4] = 0.7645655513645124
$arrDH[1,0] = 0.009659313430020688
$arrDH[1,1] = 0.08039564552838119
$arrDH[1,2] = 0.9725066926833534
$arrDH[1,3] = 0.942272679903482
$arrDH[1,4] = 0.7206173081418115
$arrDH[2,0] = 0.008846725707257974
$arrDH[2,1] = 0.08172991297474841
$arrDH[2,2] = 0.9171550594140285
$arrDH[2,3] = 0.876593798487022
$arrDH[2,4] = 0.6942582471556022
$arrDH[3,0] = 0.008517889251933752
$arrDH[3,1] = 0.08229364610825041
$arrDH[3,2] = 0.8947939911865177
$arrDH[3,3] = 0.8500019031600345
$arrDH[3,4] = 0.6836714243439133
$arrDH[4,0] = 0.00846342103194786
$arrDH[4,1] = 0.08238846180226389
$arrDH[4,2] = 0.891092602215636
$arrDH[4,3] = 0.8455966083115811
$arrDH[4,4] = 0.6819227609746576
$arrDH[5,0] = 0.008842281776360039
$arrDH[5,1] = 0.08173743467305616
$arrDH[5,2] = 0.9168527062552414
$arrDH[5,3] = 0.8762344788380005
$arrDH[5,4] = 0.6941148464769924
$arrDH[6,0] = 0.01053467589823498
$arrDH[6,1] = 0.07904027162168514
$arrDH[6,2] = 1.032252643382918
$arrDH[6,3] = 1.012964948236231
$arrDH[6,4] = 0.749280856065667
$arrDH[7,0] = 0.01393697086607659
$arrDH[7,1] = 0.07438304367712689
$arrDH[7,2] = 1.265090107237711
$arrDH[7,3] = 1.287000617929607
$arrDH[7,4] = 0.8625448244238783
$arrDH[8,0] = 0.01649687273585698
$arrDH[8,1] = 0.07134699480116957
$arrDH[8,2] = 1.440479211173226
$arrDH[8,3] = 1.492338520280271
$arrDH[8,4] = 0.9490383187026623
$arrDH[9,0] = 0.01767658116322934
$arrDH[9,1] = 0.07004964879288167
$arrDH[9,2] = 1.521282223839165
$arrDH[9,3] = 1.586710997141779
$arrDH[9,4] = 0.9891383724448701
$arrDH[10,0] = 0.01812565882629258
$arrDH[10,1] = 0.06957043605755331
$arrDH[10,2] = 1.55203244377887
$arrDH[10,3] = 1.622592937645095
$arrDH[10,4] = 1.004434751645817
$arrDH[11,0] = 0.01802883488264229
$arrDH[11,1] = 0.06967310624024137
$arrDH[11,2] = 1.545402989834372
$arrDH[11,3] = 1.614858558290223
$arrDH[11,4] = 1.001135401981685
$arrDH[12,0] = 0.01771347905939535
$arrDH[12,1] = 0.07000998178313722
$arrDH[12,2] = 1.523808987087392
$arrDH[12,3] = 1.589660078977886
$arrDH[12,4] = 0.9903945657226245
$arrDH[13,0] = 0.01752062512360197
$arrDH[13,1] = 0.07021789943411605
$arrDH[13,2] = 1.510601981363607
$arrDH[13,3] = 1.574244381763094
$arrDH[13,4] = 0.9838300848833796
$arrDH[14,0] = 0.01642009562475266
$arrDH[14,1] = 0.07143346658254046
$arrDH[14,2] = 1.435219492688844
$arrDH[14,3] = 1.486191002864757
$arrDH[14,4] = 0.9464331057195068
$arrDH[15,0] = 0.0157489632180372
$arrDH[15,1] = 0.07220064535735737
$arrDH[15,2] = 1.389239325280101
$arrDH[15,3] = 1.432424486907053
$arrDH[15,4] = 0.9236864159461788
$arrDH[16,0] = 0.01536436605358915
$arrDH[16,1] = 0.07264978977537062
$arrDH[16,2] = 1.362888195790475
$arrDH[16,3] = 1.401589670192777
$arrDH[16,4] = 0.9106738848543614
$arrDH[17,0] = 0.01523438712256819
$arrDH[17,1] = 0.07280321580784932
$arrDH[17,2] = 1.353982382094443
$arrDH[17,3] = 1.391164833777196
$arrDH[17,4] = 0.9062801293958671
$arrDH[18,0] = 0.01582025823071831
$arrDH[18,1] = 0.07211816182762831
$arrDH[18,2] = 1.394124074435268
$arrDH[18,3] = 1.438138635056646
$arrDH[18,4] = 0.9261004948243396
$arrDH[19,0] = 0.0178060416842527
$arrDH[19,1] = 0.06991070570480495
$arrDH[19,2] = 1.530147502699947
$arrDH[19,3] = 1.597057490987254
$arrDH[19,4] = 0.993546363445148
$arrDH[20,0] = 0.01911762763781866
$arrDH[20,1] = 0.06853833316785796
$arrDH[20,2] = 1.619934551999052
$arrDH[20,3] = 1.701769247391212
$arrDH[20,4] = 1.038276646083148
$arrDH[21,0] = 0.01841629586508731
$arrDH[21,1] = 0.06926435293083699
$arrDH[21,2] = 1.571930415973327
$arrDH[21,3] = 1.64580270111469
$arrDH[21,4] = 1.014342732873843
$arrDH[22,0] = 0.01578802188787876
$arrDH[22,1] = 0.07215542745558778
$arrDH[22,2] = 1.391915419734318
$arrDH[22,3] = 1.435555031614797
$arrDH[22,4] = 0.9250088880136218
$arrDH[23,0] = 0.01300666510171311
$arrDH[23,1] = 0.07557526015544447
$arrDH[23,2] = 1.201366875987105
$arrDH[23,3] = 1.212193047296438
$arrDH[23,4] = 0.831341337919838
$ws.Range("D2:H25").Value = $arrDH

$arrMN = New-Object 'object[,]' 24,2
$arrMN[0,0] = 0.9667229984368362
$arrMN[0,1] = 1.345383587036309
$arrMN[1,0] = 0.8458850692355213
$arrMN[1,1] = 1.308409149654665
$arrMN[2,0] = 0.771693125926376
$arrMN[2,1] = 1.286260518814828
$arrMN[3,0] = 0.7414586268383232
$arrMN[3,1] = 1.277375529124143
$arrMN[4,0] = 0.7364381292003088
$arrMN[4,1] = 1.275908729091128
$arrMN[5,0] = 0.771285377522787
$arrMN[5,1] = 1.286140120814153
$arrMN[6,0] = 0.92505653164676
$arrMN[6,1] = 1.332520729380633
$arrMN[7,0] = 1.226707465036256
$arrMN[7,1] = 1.427808475257649
$arrMN[8,0] = 1.44853350833678
$arrMN[8,1] = 1.500390354605315
$arrMN[9,0] = 1.549524306988843
$arrMN[9,1] = 1.533954996273195
$arrMN[10,0] = 1.587780785270041
$arrMN[10,1] = 1.546742393132575
$arrMN[11,0] = 1.579540940264593
$arrMN[11,1] = 1.543984987646382
$arrMN[12,0] = 1.552671410230403
$arrMN[12,1] = 1.535005484103294
$arrMN[13,0] = 1.536214854247476
$arrMN[13,1] = 1.529515289998784
$arrMN[14,0] = 1.441935304995027
$arrMN[14,1] = 1.49820770606965
$arrMN[15,0] = 1.38411986668325
$arrMN[15,1] = 1.479140534930906
$arrMN[16,0] = 1.350873469207826
$arrMN[16,1] = 1.468225151936082
$arrMN[17,0] = 1.339618046964461
$arrMN[17,1] = 1.464538288493685
$arrMN[18,0] = 1.390273628336004
$arrMN[18,1] = 1.481164943113555
$arrMN[19,0] = 1.560563262628904
$arrMN[19,1] = 1.537640898175027
$arrMN[20,0] = 1.671936697119889
$arrMN[20,1] = 1.5750005077509
$arrMN[21,0] = 1.61248672633721
$arrMN[21,1] = 1.555020335565416
$arrMN[22,0] = 1.387491536287854
$arrMN[22,1] = 1.480249562955152
$arrMN[23,0] = 1.145077469469825
$arrMN[23,1] = 1.401573705942781
$ws.Range("M2:N25").Value = $arrMN
